$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row: duplicate row 351 down into a new row 352,
# pushing the former rows 352..366 down to 353..367.
$ws.Rows.Item(351).Copy()
$ws.Rows.Item(352).Insert()

# Row 351 becomes the new weekly data point; update its values.
$ws.Range("D351").Value = 45041
$ws.Range("J351").Value = 260
$ws.Range("K351").Value = 18000
$ws.Range("L351").Value = 20000
$ws.Range("M351").Value = 19077
$ws.Range("P351").Value = 763
